# Generate Report for Handback
# Updates the handback timestamps (and priority) that get refreshed when the
# handback report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" column (also mirrored onto
# the de-de sheet's "Correspond Handoff Datetime" column, which shares the
# same value).
$wsOverview.Range("G2").Value = "2016-09-07 10:15:47"
$wsOverview.Range("G3").Value = "2016-09-07 10:15:47"

# zh-cn sheet - Priority changed from "ht" to "mt" (also mirrored onto the
# de-de sheet's Priority column, which shares the same value).
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet - "Correspond Handoff Datetime" column.
$wsZhCn.Range("H2").Value = "2016-09-07 10:15:35"
$wsZhCn.Range("H3").Value = "2016-09-07 10:15:35"

# zh-cn sheet - "Correspond Handback DateTime" column.
$wsZhCn.Range("K2").Value = "2016-09-07 10:16:27"
$wsZhCn.Range("K3").Value = "2016-09-07 10:16:27"

# de-de sheet - "Correspond Handoff Datetime" column (mirrors Overview G).
$wsDeDe.Range("H2").Value = "2016-09-07 10:15:47"
$wsDeDe.Range("H3").Value = "2016-09-07 10:15:47"

# de-de sheet - "Correspond Handback DateTime" column.
$wsDeDe.Range("K2").Value = "2016-09-07 10:16:44"
$wsDeDe.Range("K3").Value = "2016-09-07 10:16:44"

# de-de sheet - Priority column (mirrors zh-cn E).
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
